# Applies the commit's changes to the "dentro_rotonda_test.xlsx" workbook.
#
# Summary of the edit (per the supplied diff):
#   Caco_3 (sheet r"A1:N11"->"A1:N9"):
#     - row 9 is overwritten with what used to be the trailing "salgo rotonda"/
#       "RND-EXIT" sample, and the old trailing rows 10 ("sigo recto") and 11
#       ("acelero") are removed entirely.
#   Caco_5 (sheet "A1:N7"->"A1:N8"):
#     - a brand-new row is inserted at row 5 ("quito intermitente"/"BLK-OFF"),
#       pushing the old rows 5-7 down to 6-8 unchanged.
#   Caco_6, Caco_7, Caco_8:
#     - simple typo fixes: "gira derecha"/"gira izquierda" -> "giro derecha"/
#       "giro izquierda" (no other data changes).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Caco_3: replace row 9, then delete the old rows 10 and 11.
# ---------------------------------------------------------------------------
$wsCaco3 = $wb.Worksheets.Item("Caco_3")

$wsCaco3.Range("A9").Value = "salgo rotonda"
$wsCaco3.Range("B9").Value = "RND-EXIT"
$wsCaco3.Range("C9").Value = 4
$wsCaco3.Range("D9").Value = 0
$wsCaco3.Range("E9").Value = 0.8500000238418579
$wsCaco3.Range("F9").Value = -0.800000011920929
$wsCaco3.Range("G9").Value = 0
$wsCaco3.Range("H9").Value = "Right_Blinker"
$wsCaco3.Range("I9").Value = $false
$wsCaco3.Range("J9").Value = $false
$wsCaco3.Range("K9").Value = 34
$wsCaco3.Range("L9").Value = 11.56580257415771
$wsCaco3.Range("M9").Value = -19.83598136901855
$wsCaco3.Range("N9").Value = 0.002488384256139398

# Remove what used to be rows 10 ("sigo recto") and 11 ("acelero"). Deleting
# row index 10 twice removes both, since each delete shifts the rows below up.
$wsCaco3.Rows.Item(10).Delete()
$wsCaco3.Rows.Item(10).Delete()

# ---------------------------------------------------------------------------
# Caco_5: insert a new row 5 and populate it; old rows 5-7 shift to 6-8.
# ---------------------------------------------------------------------------
$wsCaco5 = $wb.Worksheets.Item("Caco_5")

$wsCaco5.Rows.Item(5).Insert()

$wsCaco5.Range("A5").Value = "quito intermitente"
$wsCaco5.Range("B5").Value = "BLK-OFF"
$wsCaco5.Range("C5").Value = 2
$wsCaco5.Range("D5").Value = 0
$wsCaco5.Range("E5").Value = 0.5020059943199158
$wsCaco5.Range("F5").Value = -0.2378262132406235
$wsCaco5.Range("G5").Value = 0
$wsCaco5.Range("H5").Value = "Left_Blinker"
$wsCaco5.Range("I5").Value = $false
$wsCaco5.Range("J5").Value = $false
$wsCaco5.Range("K5").Value = 41
$wsCaco5.Range("L5").Value = -20.16416549682617
$wsCaco5.Range("M5").Value = -7.858899593353271
$wsCaco5.Range("N5").Value = 0.001834316179156303

# ---------------------------------------------------------------------------
# Caco_6, Caco_7, Caco_8: plain text corrections ("gira" -> "giro").
# ---------------------------------------------------------------------------
$wsCaco6 = $wb.Worksheets.Item("Caco_6")
$wsCaco6.Range("A7").Value = "giro derecha"

$wsCaco7 = $wb.Worksheets.Item("Caco_7")
$wsCaco7.Range("A5").Value = "giro izquierda"

$wsCaco8 = $wb.Worksheets.Item("Caco_8")
$wsCaco8.Range("A8").Value = "giro derecha"
